$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D (shifts old D:K -> F:M)
$ws.Range("D1:E1").EntireColumn.Insert()

# Copy number formats/styles from column F (the old column D, now shifted)
# into the two newly inserted blank columns D:E for every data row so the
# new cells inherit the same date/number styling as their neighbours.
$ws.Range("F7:F102").Copy()
$ws.Range("D7:E7").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the two new columns with the latest quarter figures
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 31000
$ws.Range("E8").Value = 30800
$ws.Range("D9:E9").Value = "NA"
$ws.Range("D10:E10").Value = "NA"
$ws.Range("D12:E12").Value = "NA"
$ws.Range("D13:E13").Value = 0
$ws.Range("D14:E14").Value = 0
$ws.Range("D15:E15").Value = -700
$ws.Range("D17").Value = 7600
$ws.Range("E17").Value = 6200
$ws.Range("D18").Value = 23400
$ws.Range("E18").Value = 24600
$ws.Range("D20").Value = -9700
$ws.Range("E20").Value = -13900
$ws.Range("D21").Value = 15200
$ws.Range("E21").Value = 12200
$ws.Range("D22:E22").Value = 0
$ws.Range("D23").Value = 13700
$ws.Range("E23").Value = 10700
$ws.Range("D24").Value = 3500
$ws.Range("E24").Value = 2500
$ws.Range("D25:E25").Value = 0
$ws.Range("D26").Value = 10200
$ws.Range("E26").Value = 8300
$ws.Range("D27").Value = 10200
$ws.Range("E27").Value = 8300
$ws.Range("D28:E28").Value = 0
$ws.Range("D29:E29").Value = "NA"
$ws.Range("D30:E30").Value = 0
$ws.Range("D31:E31").Value = 0
$ws.Range("D32").Value = 9700
$ws.Range("E32").Value = 13900
$ws.Range("D33").Value = 10200
$ws.Range("E33").Value = 8300
$ws.Range("D34:E34").Value = 0
$ws.Range("D35").Value = 10200
$ws.Range("E35").Value = 8300
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 41500
$ws.Range("E41").Value = 45800
$ws.Range("D42").Value = 3000
$ws.Range("E42").Value = 4500
$ws.Range("D43:E43").Value = 0
$ws.Range("D44:E44").Value = 0
$ws.Range("D45:E45").Value = 0
$ws.Range("D46:E46").Value = 0
$ws.Range("D47:E47").Value = 0
$ws.Range("D48").Value = 42600
$ws.Range("E48").Value = 43100
$ws.Range("D49").Value = 110000
$ws.Range("E49").Value = 110400
$ws.Range("D50:E50").Value = 0
$ws.Range("D51:E51").Value = 0
$ws.Range("D52").Value = 4700
$ws.Range("E52").Value = 11400
$ws.Range("D53:E53").Value = 0
$ws.Range("D54").Value = 2950000
$ws.Range("E54").Value = 2931000
$ws.Range("D57").Value = 2800
$ws.Range("E57").Value = 1600
$ws.Range("D58:E58").Value = 0
$ws.Range("D59:E59").Value = 0
$ws.Range("D60:E60").Value = 0
$ws.Range("D61").Value = 38400
$ws.Range("E61").Value = 38300
$ws.Range("D62").Value = 6000
$ws.Range("E62").Value = 6100
$ws.Range("D63:E63").Value = 0
$ws.Range("D64:E64").Value = 0
$ws.Range("D65:E65").Value = 0
$ws.Range("D66").Value = 2579200
$ws.Range("E66").Value = 2570800
$ws.Range("D68:E68").Value = 0
$ws.Range("D69:E69").Value = 0
$ws.Range("D70:E70").Value = 0
$ws.Range("D71:E71").Value = 0
$ws.Range("D72").Value = 82600
$ws.Range("E72").Value = 74200
$ws.Range("D73:E73").Value = 0
$ws.Range("D74:E74").Value = 0
$ws.Range("D75:E75").Value = 0
$ws.Range("D76").Value = 370800
$ws.Range("E76").Value = 360200
$ws.Range("D77:E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = 10200
$ws.Range("E81").Value = 8300
$ws.Range("D83:E83").Value = 1500
$ws.Range("D84:E84").Value = 0
$ws.Range("D85:E85").Value = 0
$ws.Range("D86:E86").Value = 0
$ws.Range("D87:E87").Value = 0
$ws.Range("D88:E88").Value = 0
$ws.Range("D89").Value = 9700
$ws.Range("E89").Value = 13300
$ws.Range("D91").Value = -400
$ws.Range("E91").Value = -100
$ws.Range("D92:E92").Value = 0
$ws.Range("D93:E93").Value = 0
$ws.Range("D94").Value = -23400
$ws.Range("E94").Value = -20200
$ws.Range("D96:E96").Value = -1700
$ws.Range("D97:E97").Value = 0
$ws.Range("D98:E98").Value = 0
$ws.Range("D99:E99").Value = 0
$ws.Range("D100").Value = 7900
$ws.Range("E100").Value = -9200
$ws.Range("D101:E101").Value = 0
$ws.Range("D102").Value = -5800
$ws.Range("E102").Value = -16200
